# Update "想去人数" (number of people interested) values on the "展览" and
# "全部类型" sheets to match the refreshed data output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 89
$ws1.Range("F5").Value = 312
$ws1.Range("F7").Value = 146
$ws1.Range("F12").Value = 129
$ws1.Range("F13").Value = 3142
$ws1.Range("F19").Value = 558
$ws1.Range("F20").Value = 27
$ws1.Range("F21").Value = 641
$ws1.Range("F27").Value = 2303
$ws1.Range("F28").Value = 4852
$ws1.Range("F30").Value = 69
$ws1.Range("F33").Value = 257
$ws1.Range("F36").Value = 478
$ws1.Range("F38").Value = 65
$ws1.Range("F39").Value = 147
$ws1.Range("F45").Value = 23

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 89
$ws4.Range("F5").Value = 312
$ws4.Range("F7").Value = 146
$ws4.Range("F12").Value = 129
$ws4.Range("F13").Value = 3142
$ws4.Range("F20").Value = 558
$ws4.Range("F21").Value = 27
$ws4.Range("F22").Value = 641
$ws4.Range("F28").Value = 2303
$ws4.Range("F29").Value = 4852
$ws4.Range("F31").Value = 69
$ws4.Range("F34").Value = 257
$ws4.Range("F37").Value = 478
$ws4.Range("F39").Value = 65
$ws4.Range("F40").Value = 147
$ws4.Range("F46").Value = 23
